$d = $word.ActiveDocument

# 1. "moderately and only when " -> "just moderately, and only when "
$d.Content.Find.Execute("moderately and only when ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "just moderately, and only when ", 2)

# 2. "often by the " -> "often by an "
$d.Content.Find.Execute("often by the ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "often by an ", 2)

# 3. "this must also be corrected for" -> "this is regularly corrected for"
$d.Content.Find.Execute("this must also be corrected for", $true, $false, $false, $false, $false,
                         $true, 1, $false, "this is regularly corrected for", 2)

# 4. "glaze must be scraped off " -> "glaze is scraped off "
$d.Content.Find.Execute("glaze must be scraped off ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "glaze is scraped off ", 2)

# 5. "paintings must be examined under ultraviolet light, and so forth." -> "paintings is examined under ultraviolet light, and so forth."
$d.Content.Find.Execute("paintings must be examined under ultraviolet light, and so forth.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "paintings is examined under ultraviolet light, and so forth.", 2)

# 6. "likely to be " -> "probably "
$d.Content.Find.Execute("likely to be ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "probably ", 2)
